# Updates cryptos list values (price & 1h volume %) per upstream data refresh,
# including an Aave/FraxShare row-order swap (rows 44-45).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.837.90'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.34%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.622.32'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.51%  '
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.78%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.47'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.37%  '
# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.15%  '
# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.80%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.22'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.60%  '
# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.07%  '
# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.20%  '
# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.16%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.851.72'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.58%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.623.60'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.09%  '
# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.96%  '
# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.45%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.29'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.00%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.833.27'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.31%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.62'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.79%  '
# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.39%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.60'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.51%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.76%  '
# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.77%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.06'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.63%  '
# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.30%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.90'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.57%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.89'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.30%  '
# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.11%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.48'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.60%  '
# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.72%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.18'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.44%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0481'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.13%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.42'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.44%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.09'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.52%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.390.48'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.85%  '
# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.88%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +11.85%  '
# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.18%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0169'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.67%  '
# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.28%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.75%  '
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.71%  '
# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.88%  '
# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.46'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.44%  '
# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.48'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.06%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.762.36'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.59%  '
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.80%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.77'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.58%  '
# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.28%  '
# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.14%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.58'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.16%  '
